$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.611.96'
$ws.Range("E2").Value = '  +5.50%  '

$ws.Range("D3").Value = '2.045.60'
$ws.Range("E3").Value = '  +3.07%  '

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = "'251.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.87%  '

$ws.Range("E6").Value = '  +2.21%  '

$ws.Range("D7").Value = "'64.80"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +15.25%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = "'0.376"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.66%  '

$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").Value = "'59.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.18%  '

$ws.Range("D11").Value = "'0.0755"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.58%  '

$ws.Range("D12").Value = "'0.103"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.53%  '

$ws.Range("D13").Value = "'0.908"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.29%  '

$ws.Range("D14").Value = "'15.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.37%  '

$ws.Range("D15").Value = '2.345.85'
$ws.Range("E15").Value = '  +3.23%  '

$ws.Range("D16").Value = "'5.58"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.29%  '

$ws.Range("D17").Value = "'20.27"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +18.20%  '

$ws.Range("D18").Value = '2.070.77'
$ws.Range("E18").Value = '  +4.40%  '

$ws.Range("D19").Value = '37.582.11'
$ws.Range("E19").Value = '  +5.82%  '

$ws.Range("E20").Value = '  +4.52%  '

$ws.Range("D21").Value = '0.0₃0872'
$ws.Range("E21").Value = '  +4.36%  '

$ws.Range("D22").Value = "'5.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.74%  '

$ws.Range("D23").Value = "'237.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.15%  '

$ws.Range("E24").Value = '  +17.21%  '

$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").Value = "'2.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.97%  '

$ws.Range("E27").Value = '  +5.21%  '

$ws.Range("D28").Value = "'158.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.69%  '

$ws.Range("D29").Value = "'19.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.15%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = "'0.122"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.44%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = "'5.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +8.55%  '

$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").Value = "'0.113"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +25.68%  '

$ws.Range("E33").Value = '  +6.39%  '

$ws.Range("D34").Value = "'4.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +11.05%  '

$ws.Range("E35").Value = '  +4.63%  '

$ws.Range("E36").Value = '  +6.22%  '

$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("B38").Value = 'WEMIXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").Value = "'1.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.03%  '

$ws.Range("B39").Value = 'THORChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D39").Value = "'6.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +23.69%  '

$ws.Range("E40").Value = '  +15.77%  '

$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = "'1.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.82%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = "'2.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +23.00%  '

$ws.Range("E43").Value = '  +3.97%  '

$ws.Range("E44").Value = '  +4.27%  '

$ws.Range("E45").Value = '  +5.24%  '

$ws.Range("E46").Value = '  +8.79%  '

$ws.Range("D47").Value = "'16.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.94%  '

$ws.Range("D48").Value = "'94.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.64%  '

$ws.Range("D49").Value = '1.424.22'
$ws.Range("E49").Value = '  +3.65%  '

$ws.Range("E50").Value = '  +2.49%  '

$ws.Range("D51").Value = "'47.40"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.67%  '
